$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 6945525.5
$ws.Range("I32").Value = 20834020
$ws.Range("J32").Value = 1278.5
$ws.Range("K32").Value = 20834020
$ws.Range("L32").Value = 1278.5
$ws.Range("M32").Value = -20833694
$ws.Range("N32").Value = -1930.5
$ws.Range("H132").Value = 6333048
$ws.Range("I132").Value = 7356217.5
$ws.Range("K132").Value = 22068652.5
$ws.Range("M132").Value = -22066122.5
$ws.Range("H137").Value = 1386.1277
$ws.Range("I137").Value = 1065.0333
$ws.Range("J137").Value = 1952.7646
$ws.Range("K137").Value = 3195.0999
$ws.Range("L137").Value = 5858.293799999999
$ws.Range("M137").Value = -645.0999000000002
$ws.Range("N137").Value = -10958.2938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 118.71429
$ws.Range("I5").Value = 113.5
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 113.5
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = -1.5
$ws.Range("N5").Value = -374
$ws.Range("H32").Value = 16134173
$ws.Range("I32").Value = 5090.9624
$ws.Range("J32").Value = 111116540
$ws.Range("K32").Value = 5090.9624
$ws.Range("L32").Value = 111116540
$ws.Range("M32").Value = -4803.9624
$ws.Range("N32").Value = -111117114
$ws.Range("H50").Value = 441.25
$ws.Range("I50").Value = 374
$ws.Range("J50").Value = 463.66666
$ws.Range("K50").Value = 374
$ws.Range("L50").Value = 463.66666
$ws.Range("M50").Value = 340
$ws.Range("N50").Value = -1891.66666
$ws.Range("H61").Value = 2924743.5
$ws.Range("I61").Value = 3402111.8
$ws.Range("J61").Value = 863
$ws.Range("K61").Value = 3402111.8
$ws.Range("L61").Value = 863
$ws.Range("M61").Value = -3401899.8
$ws.Range("N61").Value = -1287
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 40000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 40000
$ws.Range("M82").Value = $null
$ws.Range("N82").Value = -40722
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 40000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 40000
$ws.Range("M85").Value = $null
$ws.Range("N85").Value = -42496
$ws.Range("H94").Value = 35975
$ws.Range("J94").Value = 35975
$ws.Range("L94").Value = 35975
$ws.Range("N94").Value = -37777
$ws.Range("H122").Value = 1786.1111
$ws.Range("I122").Value = 1766.909
$ws.Range("J122").Value = 1816.2858
$ws.Range("K122").Value = 5300.727000000001
$ws.Range("L122").Value = 5448.857400000001
$ws.Range("M122").Value = -2850.727000000001
$ws.Range("N122").Value = -10348.8574
$ws.Range("H136").Value = 2924743.5
$ws.Range("I136").Value = 3402111.8
$ws.Range("J136").Value = 863
$ws.Range("K136").Value = 10206335.4
$ws.Range("L136").Value = 2589
$ws.Range("M136").Value = -10203785.4
$ws.Range("N136").Value = -7689

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 118.71429
$ws.Range("I4").Value = 113.5
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 113.5
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = 1.5
$ws.Range("N4").Value = -380
$ws.Range("H80").Value = 3961.0588
$ws.Range("I80").Value = 744.44446
$ws.Range("J80").Value = 5119.04
$ws.Range("K80").Value = 744.44446
$ws.Range("L80").Value = 5119.04
$ws.Range("M80").Value = 253.55554
$ws.Range("N80").Value = -7115.04
$ws.Range("H83").Value = 3961.0588
$ws.Range("I83").Value = 744.44446
$ws.Range("J83").Value = 5119.04
$ws.Range("K83").Value = 3722.2223
$ws.Range("L83").Value = 25595.2
$ws.Range("M83").Value = 1269.7777
$ws.Range("N83").Value = -35579.2
$ws.Range("H105").Value = 43479816
$ws.Range("I105").Value = 1399.4706
$ws.Range("K105").Value = 1399.4706
$ws.Range("M105").Value = 347.5293999999999
$ws.Range("H134").Value = 3925.7715
$ws.Range("I134").Value = 817.3103599999999
$ws.Range("J134").Value = 18950
$ws.Range("K134").Value = 2451.93108
$ws.Range("L134").Value = 56850
$ws.Range("M134").Value = 83.06892000000016
$ws.Range("N134").Value = -61920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 6154.2354
$ws.Range("I7").Value = 266
$ws.Range("J7").Value = 11388.223
$ws.Range("K7").Value = 266
$ws.Range("L7").Value = 11388.223
$ws.Range("M7").Value = -153
$ws.Range("N7").Value = -11614.223
$ws.Range("H35").Value = 1779.3
$ws.Range("I35").Value = 745.625
$ws.Range("J35").Value = 5914
$ws.Range("K35").Value = 745.625
$ws.Range("L35").Value = 5914
$ws.Range("M35").Value = -451.625
$ws.Range("N35").Value = -6502
$ws.Range("H58").Value = 12048690
$ws.Range("I58").Value = 23256210
$ws.Range("J58").Value = 604.15
$ws.Range("K58").Value = 23256210
$ws.Range("L58").Value = 604.15
$ws.Range("M58").Value = -23256007
$ws.Range("N58").Value = -1010.15
$ws.Range("H99").Value = 62503490
$ws.Range("I99").Value = 125002490
$ws.Range("J99").Value = 4487.5
$ws.Range("K99").Value = 125002490
$ws.Range("L99").Value = 4487.5
$ws.Range("M99").Value = -125000992
$ws.Range("N99").Value = -7483.5
$ws.Range("H105").Value = 4499.085
$ws.Range("I105").Value = 4446.772
$ws.Range("J105").Value = 5990
$ws.Range("K105").Value = 4446.772
$ws.Range("L105").Value = 5990
$ws.Range("M105").Value = -2699.772
$ws.Range("N105").Value = -9484
$ws.Range("H126").Value = 62503490
$ws.Range("I126").Value = 125002490
$ws.Range("J126").Value = 4487.5
$ws.Range("K126").Value = 375007470
$ws.Range("L126").Value = 13462.5
$ws.Range("M126").Value = -375005000
$ws.Range("N126").Value = -18402.5
$ws.Range("H134").Value = 8772675
$ws.Range("I134").Value = 923.2857
$ws.Range("J134").Value = 17241952
$ws.Range("K134").Value = 2769.8571
$ws.Range("L134").Value = 51725856
$ws.Range("M134").Value = -234.8571000000002
$ws.Range("N134").Value = -51730926
$ws.Range("H136").Value = 12048690
$ws.Range("I136").Value = 23256210
$ws.Range("J136").Value = 604.15
$ws.Range("K136").Value = 69768630
$ws.Range("L136").Value = 1812.45
$ws.Range("M136").Value = -69766080
$ws.Range("N136").Value = -6912.45

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 143140.14
$ws.Range("I4").Value = 200276.2
$ws.Range("K4").Value = 600828.6000000001
$ws.Range("M4").Value = -600716.6000000001
$ws.Range("H23").Value = 171
$ws.Range("I23").Value = 101
$ws.Range("J23").Value = 179.75
$ws.Range("K23").Value = 303
$ws.Range("L23").Value = 539.25
$ws.Range("M23").Value = -68
$ws.Range("N23").Value = -1009.25
$ws.Range("H25").Value = 639.6
$ws.Range("I25").Value = 99
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 297
$ws.Range("L25").Value = 3000
$ws.Range("N25").Value = -3338
$ws.Range("H30").Value = 639.6
$ws.Range("I30").Value = 99
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 297
$ws.Range("L30").Value = 3000
$ws.Range("N30").Value = -3204
$ws.Range("H87").Value = 9642.666999999999
$ws.Range("I87").Value = 9642.666999999999
$ws.Range("K87").Value = 28928.001
$ws.Range("M87").Value = -27680.001
$ws.Range("H90").Value = 9642.666999999999
$ws.Range("I90").Value = 9642.666999999999
$ws.Range("K90").Value = 86784.003
$ws.Range("M90").Value = -80544.003
$ws.Range("H131").Value = 915.61
$ws.Range("I131").Value = 487.5
$ws.Range("J131").Value = 933.44794
$ws.Range("K131").Value = 1462.5
$ws.Range("L131").Value = 2800.34382
$ws.Range("M131").Value = 3577.5
$ws.Range("N131").Value = -12880.34382

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1782.75
$ws.Range("I102").Value = 1815.5
$ws.Range("J102").Value = 1750
$ws.Range("K102").Value = 1815.5
$ws.Range("L102").Value = 1750
$ws.Range("M102").Value = -193.5
$ws.Range("N102").Value = -4994
$ws.Range("H113").Value = 1296.5264
$ws.Range("I113").Value = 1076.5
$ws.Range("J113").Value = 1398.0769
$ws.Range("K113").Value = 1076.5
$ws.Range("L113").Value = 1398.0769
$ws.Range("M113").Value = 1093.5
$ws.Range("N113").Value = -5738.0769
$ws.Range("H122").Value = 8477734
$ws.Range("I122").Value = 11114614
$ws.Range("K122").Value = 33343842
$ws.Range("M122").Value = -33341392

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8179.9375
$ws.Range("I122").Value = 11937.9
$ws.Range("J122").Value = 1916.6666
$ws.Range("K122").Value = 35813.7
$ws.Range("L122").Value = 5749.9998
$ws.Range("M122").Value = -33363.7
$ws.Range("N122").Value = -10649.9998
$ws.Range("H136").Value = 61510460
$ws.Range("I136").Value = 26458006
$ws.Range("J136").Value = 166667820
$ws.Range("K136").Value = 79374018
$ws.Range("L136").Value = 500003460
$ws.Range("M136").Value = -79371468
$ws.Range("N136").Value = -500008560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1408.7
$ws.Range("I126").Value = 1020.5
$ws.Range("J126").Value = 1667.5
$ws.Range("K126").Value = 3061.5
$ws.Range("L126").Value = 5002.5
$ws.Range("M126").Value = -591.5
$ws.Range("N126").Value = -9942.5
$ws.Range("H136").Value = 17859538
$ws.Range("I136").Value = 41668428
$ws.Range("J136").Value = 2869.375
$ws.Range("K136").Value = 125005284
$ws.Range("L136").Value = 8608.125
$ws.Range("M136").Value = -125002734
$ws.Range("N136").Value = -13708.125
